$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="57.721.99"'
$ws.Range("E2").Formula = '="  -4.20%  "'
$ws.Range("D2:E2").Copy() | Out-Null
$ws.Range("D2:E2").PasteSpecial(-4163) | Out-Null

$ws.Range("D3").Formula = '="3.096.39"'
$ws.Range("E3").Formula = '="  -6.07%  "'
$ws.Range("D3:E3").Copy() | Out-Null
$ws.Range("D3:E3").PasteSpecial(-4163) | Out-Null

$ws.Range("E4").Formula = '="  +0.14%  "'
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4163) | Out-Null

$ws.Range("D5").Formula = '="517.61"'
$ws.Range("E5").Formula = '="  -7.18%  "'
$ws.Range("D5:E5").Copy() | Out-Null
$ws.Range("D5:E5").PasteSpecial(-4163) | Out-Null

$ws.Range("D6").Formula = '="130.30"'
$ws.Range("E6").Formula = '="  -7.44%  "'
$ws.Range("D6:E6").Copy() | Out-Null
$ws.Range("D6:E6").PasteSpecial(-4163) | Out-Null

$ws.Range("E7").Formula = '="  -0.03%  "'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4163) | Out-Null

$ws.Range("D8").Formula = '="3.094.00"'
$ws.Range("E8").Formula = '="  -6.19%  "'
$ws.Range("D8:E8").Copy() | Out-Null
$ws.Range("D8:E8").PasteSpecial(-4163) | Out-Null

$ws.Range("D9").Formula = '="0.438"'
$ws.Range("E9").Formula = '="  -6.24%  "'
$ws.Range("D9:E9").Copy() | Out-Null
$ws.Range("D9:E9").PasteSpecial(-4163) | Out-Null

$ws.Range("D10").Formula = '="7.19"'
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null

$ws.Range("E11").Formula = '="  -10.68%  "'
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4163) | Out-Null

$ws.Range("D12").Formula = '="0.367"'
$ws.Range("E12").Formula = '="  -9.66%  "'
$ws.Range("D12:E12").Copy() | Out-Null
$ws.Range("D12:E12").PasteSpecial(-4163) | Out-Null

$ws.Range("D13").Formula = '="3.634.19"'
$ws.Range("E13").Formula = '="  -5.88%  "'
$ws.Range("D13:E13").Copy() | Out-Null
$ws.Range("D13:E13").PasteSpecial(-4163) | Out-Null

$ws.Range("E14").Formula = '="  -0.91%  "'
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4163) | Out-Null

$ws.Range("D15").Formula = '="24.66"'
$ws.Range("E15").Formula = '="  -6.97%  "'
$ws.Range("D15:E15").Copy() | Out-Null
$ws.Range("D15:E15").PasteSpecial(-4163) | Out-Null

$ws.Range("D16").Formula = '="57.820.31"'
$ws.Range("E16").Formula = '="  -4.02%  "'
$ws.Range("D16:E16").Copy() | Out-Null
$ws.Range("D16:E16").PasteSpecial(-4163) | Out-Null

$ws.Range("D17").Formula = '="3.101.70"'
$ws.Range("E17").Formula = '="  -5.99%  "'
$ws.Range("D17:E17").Copy() | Out-Null
$ws.Range("D17:E17").PasteSpecial(-4163) | Out-Null

$ws.Range("D18").Formula = '="0.0000148"'
$ws.Range("E18").Formula = '="  -9.64%  "'
$ws.Range("D18:E18").Copy() | Out-Null
$ws.Range("D18:E18").PasteSpecial(-4163) | Out-Null

$ws.Range("D19").Formula = '="5.62"'
$ws.Range("E19").Formula = '="  -7.31%  "'
$ws.Range("D19:E19").Copy() | Out-Null
$ws.Range("D19:E19").PasteSpecial(-4163) | Out-Null

$ws.Range("D20").Formula = '="12.74"'
$ws.Range("E20").Formula = '="  -6.62%  "'
$ws.Range("D20:E20").Copy() | Out-Null
$ws.Range("D20:E20").PasteSpecial(-4163) | Out-Null

$ws.Range("D21").Formula = '="7.72"'
$ws.Range("E21").Formula = '="  -9.47%  "'
$ws.Range("D21:E21").Copy() | Out-Null
$ws.Range("D21:E21").PasteSpecial(-4163) | Out-Null

$ws.Range("D22").Formula = '="335.66"'
$ws.Range("E22").Formula = '="  -10.08%  "'
$ws.Range("D22:E22").Copy() | Out-Null
$ws.Range("D22:E22").PasteSpecial(-4163) | Out-Null

$ws.Range("E23").Formula = '="  -0.09%  "'
$ws.Range("E23").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4163) | Out-Null

$ws.Range("D24").Formula = '="0.501"'
$ws.Range("E24").Formula = '="  -5.65%  "'
$ws.Range("D24:E24").Copy() | Out-Null
$ws.Range("D24:E24").PasteSpecial(-4163) | Out-Null

$ws.Range("D25").Formula = '="66.14"'
$ws.Range("E25").Formula = '="  -8.13%  "'
$ws.Range("D25:E25").Copy() | Out-Null
$ws.Range("D25:E25").PasteSpecial(-4163) | Out-Null

$ws.Range("E26").Formula = '="  -3.95%  "'
$ws.Range("E26").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4163) | Out-Null

$ws.Range("E27").Formula = '="  +0.17%  "'
$ws.Range("E27").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4163) | Out-Null

$ws.Range("D28").Formula = '="0.0₃0903"'
$ws.Range("E28").Formula = '="  -11.52%  "'
$ws.Range("D28:E28").Copy() | Out-Null
$ws.Range("D28:E28").PasteSpecial(-4163) | Out-Null

$ws.Range("E29").Formula = '="  +0.01%  "'
$ws.Range("E29").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4163) | Out-Null

$ws.Range("D30").Formula = '="6.66"'
$ws.Range("E30").Formula = '="  -5.28%  "'
$ws.Range("D30:E30").Copy() | Out-Null
$ws.Range("D30:E30").PasteSpecial(-4163) | Out-Null

$ws.Range("D31").Formula = '="1.24"'
$ws.Range("E31").Formula = '="  -0.98%  "'
$ws.Range("D31:E31").Copy() | Out-Null
$ws.Range("D31:E31").PasteSpecial(-4163) | Out-Null

$ws.Range("D32").Formula = '="1.83"'
$ws.Range("E32").Formula = '="  -9.38%  "'
$ws.Range("D32:E32").Copy() | Out-Null
$ws.Range("D32:E32").PasteSpecial(-4163) | Out-Null

$ws.Range("D33").Formula = '="6.74"'
$ws.Range("E33").Formula = '="  -7.50%  "'
$ws.Range("D33:E33").Copy() | Out-Null
$ws.Range("D33:E33").PasteSpecial(-4163) | Out-Null

$ws.Range("D34").Formula = '="20.98"'
$ws.Range("E34").Formula = '="  -6.98%  "'
$ws.Range("D34:E34").Copy() | Out-Null
$ws.Range("D34:E34").PasteSpecial(-4163) | Out-Null

$ws.Range("D35").Formula = '="157.66"'
$ws.Range("E35").Formula = '="  -4.72%  "'
$ws.Range("D35:E35").Copy() | Out-Null
$ws.Range("D35:E35").PasteSpecial(-4163) | Out-Null

$ws.Range("D36").Formula = '="4.70"'
$ws.Range("E36").Formula = '="  -6.57%  "'
$ws.Range("D36:E36").Copy() | Out-Null
$ws.Range("D36:E36").PasteSpecial(-4163) | Out-Null

$ws.Range("D37").Formula = '="6.04"'
$ws.Range("E37").Formula = '="  -8.79%  "'
$ws.Range("D37:E37").Copy() | Out-Null
$ws.Range("D37:E37").PasteSpecial(-4163) | Out-Null

$ws.Range("D38").Formula = '="1.34"'
$ws.Range("E38").Formula = '="  -11.65%  "'
$ws.Range("D38:E38").Copy() | Out-Null
$ws.Range("D38:E38").PasteSpecial(-4163) | Out-Null

$ws.Range("D39").Formula = '="3.130.60"'
$ws.Range("E39").Formula = '="  -5.87%  "'
$ws.Range("D39:E39").Copy() | Out-Null
$ws.Range("D39:E39").PasteSpecial(-4163) | Out-Null

$ws.Range("D40").Formula = '="40.08"'
$ws.Range("E40").Formula = '="  -4.05%  "'
$ws.Range("D40:E40").Copy() | Out-Null
$ws.Range("D40:E40").PasteSpecial(-4163) | Out-Null

$ws.Range("D41").Formula = '="0.0669"'
$ws.Range("E41").Formula = '="  -7.20%  "'
$ws.Range("D41:E41").Copy() | Out-Null
$ws.Range("D41:E41").PasteSpecial(-4163) | Out-Null

$ws.Range("D42").Formula = '="22.90"'
$ws.Range("E42").Formula = '="  -10.58%  "'
$ws.Range("D42:E42").Copy() | Out-Null
$ws.Range("D42:E42").PasteSpecial(-4163) | Out-Null

$ws.Range("D43").Formula = '="0.680"'
$ws.Range("E43").Formula = '="  -8.71%  "'
$ws.Range("D43:E43").Copy() | Out-Null
$ws.Range("D43:E43").PasteSpecial(-4163) | Out-Null

$ws.Range("D44").Formula = '="3.85"'
$ws.Range("E44").Formula = '="  -5.82%  "'
$ws.Range("D44:E44").Copy() | Out-Null
$ws.Range("D44:E44").PasteSpecial(-4163) | Out-Null

$ws.Range("B45").Formula = '="FirstDigitalUSD"'
$ws.Range("C45").Formula = '="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"'
$ws.Range("D45").Formula = '="1.00"'
$ws.Range("E45").Formula = '="  -0.03%  "'
$ws.Range("B45:E45").Copy() | Out-Null
$ws.Range("B45:E45").PasteSpecial(-4163) | Out-Null

$ws.Range("B46").Formula = '="ONDO"'
$ws.Range("C46").Formula = '="https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"'
$ws.Range("D46").Formula = '="1.05"'
$ws.Range("E46").Formula = '="  -6.05%  "'
$ws.Range("B46:E46").Copy() | Out-Null
$ws.Range("B46:E46").PasteSpecial(-4163) | Out-Null

$ws.Range("D47").Formula = '="2.249.38"'
$ws.Range("E47").Formula = '="  -3.18%  "'
$ws.Range("D47:E47").Copy() | Out-Null
$ws.Range("D47:E47").PasteSpecial(-4163) | Out-Null

$ws.Range("D48").Formula = '="1.40"'
$ws.Range("E48").Formula = '="  -10.25%  "'
$ws.Range("D48:E48").Copy() | Out-Null
$ws.Range("D48:E48").PasteSpecial(-4163) | Out-Null

$ws.Range("D49").Formula = '="6.06"'
$ws.Range("E49").Formula = '="  -4.88%  "'
$ws.Range("D49:E49").Copy() | Out-Null
$ws.Range("D49:E49").PasteSpecial(-4163) | Out-Null

$ws.Range("D50").Formula = '="20.06"'
$ws.Range("E50").Formula = '="  -6.56%  "'
$ws.Range("D50:E50").Copy() | Out-Null
$ws.Range("D50:E50").PasteSpecial(-4163) | Out-Null

$ws.Range("D51").Formula = '="0.0229"'
$ws.Range("E51").Formula = '="  -8.92%  "'
$ws.Range("D51:E51").Copy() | Out-Null
$ws.Range("D51:E51").PasteSpecial(-4163) | Out-Null
